$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of work-hour entries
$ws.Range("B6").Value = 0.375
$ws.Range("C6").Value = 0.3888888888888889
$ws.Range("D6").Value = "Dashboard modelling with Probuilder"

$ws.Range("B7").Value = 0.3923611111111111
$ws.Range("C7").Value = 0.41666666666666669
$ws.Range("D7").Value = "Pedals and gearbox"

$ws.Range("B8").Value = 0.43055555555555558
$ws.Range("C8").Value = 0.46875
$ws.Range("D8").Value = "Pedals and gearbox"

# Match number formatting used by existing time cells (B3:C5)
$ws.Range("B6:C8").NumberFormat = $ws.Range("B5:C5").NumberFormat

# Update selection to match the final state
$ws.Range("D8").Select()
